$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 5.590999999999999
$ws.Range("C4").Value = -13.038
$ws.Range("D4").Value = -8.007999999999999

$ws.Range("C5").Value = -12.413

$ws.Range("B7").Value = 6.957000000000001

$ws.Range("C8").Value = -12.8

$ws.Range("D9").Value = -7.946

$ws.Range("B16").Value = 6.107
$ws.Range("C16").Value = -12.449

$ws.Range("D18").Value = -7.835999999999999
